# ProposalForm (sheet4) gains two new trailing columns — AQ "nomShare" and
# AR "sumAssured" — and the existing nominee DOB year (AN2) is corrected
# from 1999 to 1974. Order of writes matters because new shared-string
# table entries are appended in write order, and the target workbook's
# sharedStrings.xml ends with: nomShare, 100, 1974, sumAssured.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProposalForm")

# --- New column AQ: nomShare ------------------------------------------------
$ws.Range("AQ1").Value = "nomShare"
$ws.Range("AQ1").Font.Bold = $true

# stored as text (matches the quote-prefixed "100" already used elsewhere,
# e.g. F2/G2/AL2/AM2) rather than a genuine number
$ws.Range("AQ2").Value = "'100"

# --- Fix existing nominee DOB year (AN2): 1999 -> 1974 ----------------------
# keep it text (quote-prefixed), same as the original cell
$ws.Range("AN2").Value = "'1974"

# --- New column AR: sumAssured ----------------------------------------------
$ws.Range("AR1").Value = "sumAssured"
$ws.Range("AR1").Font.Bold = $true

# numeric value, no special formatting
$ws.Range("AR2").Value = 500000

# --- Column width for the new AQ/AR columns (best-fit-like) ----------------
$ws.Columns.Item(42).ColumnWidth = 15

# --- View / selection -------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 37
$excel.ActiveWindow.ScrollRow = 1
[void]$ws.Range("AR2").Select()
